$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.938.80"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.648.42"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5109"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2583"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06428"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.74"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07778"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.321"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "1.654.45"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5483"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "0.0₅7902"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.97"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "26.024.91"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.40"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.438"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.072"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.857"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1148"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.920"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.78"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.242"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05027"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.290"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.211"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.546"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.372"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8966"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.135.95"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5546"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01567"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.009"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.675"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8176"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.99"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("E44").Value = "  +8.17%  "
$ws.Range("D45").Value = "1.785.35"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4532"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05096"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.009"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09587"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.08%  "
